$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update D8 value (E8/F8 formulas recalc automatically)
$ws.Range("D8").Value = 7961

# 2. Copy formatting from row 10 down to new rows 11-20
$ws.Range("C10:F10").Copy()
$ws.Range("C11:F20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. Row 11 + 12 (shared formulas across E11:E12 and F11:F12)
$ws.Range("C11").Value = "FFT (DSPF_sp_cfftr2_dit)"
$ws.Range("D11").Value = 4839
$ws.Range("D12").Value = 98103
$ws.Range("E11:E12").Formula = "=D11/225000000*1000"
$ws.Range("F11:F12").Formula = "=E11/10"

# 4. Row 13
$ws.Range("C13").Value = "freq2mel"
$ws.Range("D13").Value = 1646
$ws.Range("E13").Formula = "=D13/225000000*1000"
$ws.Range("F13").Formula = "=E13/10"

# 5. Row 14
$ws.Range("C14").Value = "mel2freq"
$ws.Range("D14").Value = 1278
$ws.Range("E14").Formula = "=D14/225000000*1000"
$ws.Range("F14").Formula = "=E14/10"

# 6. Row 15
$ws.Range("C15").Value = "MvgAvg"
$ws.Range("D15").Value = 253
$ws.Range("E15").Formula = "=D15/225000000*1000"
$ws.Range("F15").Formula = "=E15/10"

# 7. Row 16
$ws.Range("C16").Value = "MelFilterBank Create"
$ws.Range("D16").Value = 323319
$ws.Range("E16").Formula = "=D16/225000000*1000"
$ws.Range("F16").Formula = "=E16/10"

# 8. Row 17
$ws.Range("C17").Value = "get MelCoeff"
$ws.Range("D17").Value = 12225
$ws.Range("E17").Formula = "=D17/225000000*1000"
$ws.Range("F17").Formula = "=E17/10"

# 9. Row 18
$ws.Range("C18").Value = "transfert between buffer 256"
$ws.Range("D18").Value = 6154
$ws.Range("E18").Formula = "=D18/225000000*1000"
$ws.Range("F18").Formula = "=E18/10"

# 10. Row 19
$ws.Range("C19").Value = "float2complex"
$ws.Range("D19").Value = 16169
$ws.Range("E19").Formula = "=D19/225000000*1000"
$ws.Range("F19").Formula = "=E19/10"

# 11. Row 12 label (added after row 19's text so shared-string index lands at 15)
$ws.Range("C12").Value = "FFT (DSPF_sp_cfftr2_dit) init"

# 12. Row 20
$ws.Range("C20").Value = "autocorrelation 256"
$ws.Range("D20").Value = 203980
$ws.Range("E20").Formula = "=D20/225000000*1000"
$ws.Range("F20").Formula = "=E20/10"

# 13. Row 21: empty styled cells E21, F21 only
$ws.Range("E10:F10").Copy()
$ws.Range("E21:F21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 14. Row 32: F32 total. Set the formula first, then apply formatting
#     afterwards to avoid a stale cached value from the paste operation.
$ws.Range("F32").Formula = "=SUM(F8:F31)"
$ws.Range("F20").Copy()
$ws.Range("F32").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 15. Update sheet view to match target (scrolled to row 7, selection on D23)
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("D23").Select() | Out-Null
